# Kernel Dev Guide overview.pptx - remove pre-installed applications row
# and relabel/resize the "downloaded applications" row to
# "Sandboxed Applications", plus bump the footer date.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Footer date placeholder: "Apr. 17" -> "Dec" + ". 2024" (two runs,
#    matching how PowerPoint splits a retyped/partial edit into runs).
# ---------------------------------------------------------------------
$dateShape = $s.Shapes.Item("Date Placeholder 2")
$dateRange = $dateShape.TextFrame.TextRange
$dateRange.Text = "Dec"
$dateRange.InsertAfter(". 2024") | Out-Null

# ---------------------------------------------------------------------
# 2) Resize/reposition the three stacked "Downloaded Applications"
#    rounded rectangles on the right-hand stack so they stretch across
#    the full width (the left-hand stack's boxes are being removed).
#    Point values below are chosen so the saved EMU lands exactly on
#    the target offsets/extents.
# ---------------------------------------------------------------------
$rect35 = $s.Shapes.Item("Rounded Rectangle 35")
$rect35.Left = 150.21291358582675
$rect35.Top = 99.90236290472441
$rect35.Width = 418.0832977866142
$rect35.Height = 53.80409438818897

$rect36 = $s.Shapes.Item("Rounded Rectangle 36")
$rect36.Left = 145.08646399291337
$rect36.Top = 105.448032396063
$rect36.Width = 418.0832977866142
$rect36.Height = 53.80409438818897

$rect37 = $s.Shapes.Item("Rounded Rectangle 37")
$rect37.Left = 139.8633041866142
$rect37.Top = 111.44590381181102
$rect37.Width = 418.0832977866142
$rect37.Height = 53.80409438818897

# ---------------------------------------------------------------------
# 3) The top-most box of that stack ("Rounded Rectangle 37") becomes a
#    single line of text "Sandboxed Applications" instead of the two
#    "Downloaded" / "Applications" paragraphs.
# ---------------------------------------------------------------------
$frontRange = $rect37.TextFrame.TextRange
$frontRange.Text = "Sandboxed Applications"
$frontRange.Font.Name = "Calibri Light"
$frontRange.Font.Color.RGB = 5722955   # 0x57534B -> srgbClr 4B5357

# ---------------------------------------------------------------------
# 4) Drop the left-hand stack entirely (the "pre-installed
#    applications" boxes): Rounded Rectangle 41/42/43.
# ---------------------------------------------------------------------
$s.Shapes.Item("Rounded Rectangle 41").Delete()
$s.Shapes.Item("Rounded Rectangle 42").Delete()
$s.Shapes.Item("Rounded Rectangle 43").Delete()
